# Indicator 5.5.1 workbook update
# - relabel the national-parliament row so it doesn't collide with the new
#   local-government row
# - correct the 2010 (column G) value
# - add a 2022 column (S)
# - add a new row for "Proportion of seats held by women in local government"
# - add a footnote row explaining the local-government source
# - move the saved selection cursor

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Values
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "Улуттук парламенттердеги аялдардын орундарынын үлүшү"
$ws.Range("G5").Value = 23.9

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 21.1

$ws.Range("A6").Value = "Жергиликтүү өз алдынча башкаруу органдарындагы аялдардын орундарынын үлүшү*"
$ws.Range("B6").Value = "Доля мест, занимаемых женщинами в местных органах власти*"
$ws.Range("C6").Value = "Proportion of seats held by women in local government*"
$ws.Range("M6").Value = 15.61
$ws.Range("N6").Value = 15.09
$ws.Range("O6").Value = 14.96
$ws.Range("P6").Value = 15.16
$ws.Range("Q6").Value = 14.98
$ws.Range("R6").Value = 31.55
$ws.Range("S6").Value = 36.46

$ws.Range("A7").Value = "*КР ШРӨБК маалыматтары боюнча"
$ws.Range("B7").Value = "*по данным ЦКВПР КР"
$ws.Range("C7").Value = "*according to the CCER of KR"

# ---------------------------------------------------------------------
# 2) Formats: extend the existing look across the new column/rows
# ---------------------------------------------------------------------
# Column S (2022) takes on the same per-cell look as column R
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("R4:S4").NumberFormat = "General"

# Row 6 (local-government row) starts from row 5's look
$ws.Range("A5:S5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("D6:S6").NumberFormat = "0.0"
$ws.Range("D6:L6").ClearContents()

# Row 7 (footnote) - small plain Times New Roman text, no border/bold/wrap
$r7 = $ws.Range("A7:C7")
$r7.Font.Name = "Times New Roman"
$r7.Font.Size = 9
$r7.Font.Bold = $false
$r7.Font.Italic = $false
$r7.Borders.Item(9).LineStyle = -4142
$r7.WrapText = $false
$r7.VerticalAlignment = -4160

# ---------------------------------------------------------------------
# 3) Borders: the table's thick rule moves from the bottom of row 5 to
#    the bottom of the new last row (row 6)
# ---------------------------------------------------------------------
$ws.Range("A5:S5").Borders.Item(9).LineStyle = -4142
$ws.Range("A6:S6").Borders.Item(9).Weight = -4138

# ---------------------------------------------------------------------
# 4) Row heights
# ---------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 27
$ws.Rows.Item(6).RowHeight = 28.5

# ---------------------------------------------------------------------
# 5) View state
# ---------------------------------------------------------------------
$ws.Range("T4").Select()

Write-Host "edit applied"
